$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample row (id 49 / row 50) — a phishing attempt targeting MCAST
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "mail"
$ws.Range("C50").Value = 44390
$ws.Range("C50").NumberFormat = "mm-dd-yy"
$ws.Range("D50").Value = "MCAST"
$ws.Range("E50").Value = "other"
$ws.Range("F50").Value = "lockout"
$ws.Range("G50").Value = "en"
$ws.Range("H50").Value = "no"
$ws.Range("I50").Value = "Verify to avoid lockout"
$ws.Range("J50").Value = "MCAST"
$ws.Range("K50").Value = "redirects to https://officehotmail2021.weebly.com/"

$ws.Range("A51").Select() | Out-Null
